# incorporate v3 data into toi data table
$wb = $excel.ActiveWorkbook

# Sheet: ColumnHeadersToi -- add new attribute row for iode_quality_flag
$wsToi = $wb.Worksheets.Item("ColumnHeadersToi")
$wsToi.Range("A17").Value = "iode_quality_flag"
$wsToi.Range("B17").Value = "IODE Quality Flag primary level"
$wsToi.Range("C17").Value = "categorical"

# Sheet: CategoricalVariables -- add the code/definition rows for iode_quality_flag
$wsCat = $wb.Worksheets.Item("CategoricalVariables")
$wsCat.Range("A4").Value = "iode_quality_flag"
$wsCat.Range("B4").Value = 1
$wsCat.Range("C4").Value = "good"

$wsCat.Range("A5").Value = "iode_quality_flag"
$wsCat.Range("B5").Value = 2
$wsCat.Range("C5").Value = "quality not evaluated, not available or unknown"

$wsCat.Range("A6").Value = "iode_quality_flag"
$wsCat.Range("B6").Value = 3
$wsCat.Range("C6").Value = "questionable/suspect"

$wsCat.Range("A7").Value = "iode_quality_flag"
$wsCat.Range("B7").Value = 4
$wsCat.Range("C7").Value = "bad"

$wsCat.Range("A8").Value = "iode_quality_flag"
$wsCat.Range("B8").Value = 9
$wsCat.Range("C8").Value = "missing data"

# Update selections/active sheet to match final state
$wsToi.Range("A17:C17").Select() | Out-Null
$wsCat.Range("A4:C8").Select() | Out-Null
$wsCat.Activate() | Out-Null
